# Refresh cryptos.xlsx price/volume figures and swap the Aptos/Algorand row contents
# (generated to match the "Updated cryptos list ... with GitHub Actions" commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.022.24"
$ws.Range("E2").Value = "'  -3.16%  "
$ws.Range("D3").Value = "'1.714.63"
$ws.Range("E3").Value = "'  -2.96%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'310.05"
$ws.Range("E5").Value = "'  -5.63%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  +0.05%  "
$ws.Range("D7").Value = "'0.4743"
$ws.Range("E7").Value = "'  +4.74%  "
$ws.Range("D8").Value = "'0.3449"
$ws.Range("E8").Value = "'  -2.07%  "
$ws.Range("D9").Value = "'42.15"
$ws.Range("E9").Value = "'  +0.53%  "
$ws.Range("D10").Value = "'0.07239"
$ws.Range("E10").Value = "'  -1.86%  "
$ws.Range("D11").Value = "'1.036"
$ws.Range("E11").Value = "'  -5.18%  "
$ws.Range("E12").Value = "'  +0.04%  "
$ws.Range("D13").Value = "'19.74"
$ws.Range("E13").Value = "'  -4.62%  "
$ws.Range("D14").Value = "'5.817"
$ws.Range("E14").Value = "'  -3.11%  "
$ws.Range("D15").Value = "'1.723.01"
$ws.Range("E15").Value = "'  -2.70%  "
$ws.Range("D16").Value = "'6.810"
$ws.Range("E16").Value = "'  -5.13%  "
$ws.Range("D17").Value = "'86.91"
$ws.Range("E17").Value = "'  -6.01%  "
$ws.Range("E18").Value = "'  -2.40%  "
$ws.Range("D19").Value = "'0.06370"
$ws.Range("E19").Value = "'  -1.13%  "
$ws.Range("E21").Value = "'  -3.01%  "
$ws.Range("D22").Value = "'5.604"
$ws.Range("D23").Value = "'27.076.26"
$ws.Range("E23").Value = "'  -3.05%  "
$ws.Range("D24").Value = "'10.72"
$ws.Range("E24").Value = "'  -4.18%  "
$ws.Range("D25").Value = "'2.099"
$ws.Range("E25").Value = "'  -0.07%  "
$ws.Range("D26").Value = "'19.92"
$ws.Range("E26").Value = "'  -1.05%  "
$ws.Range("D27").Value = "'150.63"
$ws.Range("E27").Value = "'  -4.96%  "
$ws.Range("D28").Value = "'1.917.81"
$ws.Range("E28").Value = "'  -2.82%  "
$ws.Range("D29").Value = "'2.060"
$ws.Range("E29").Value = "'  -3.76%  "
$ws.Range("D30").Value = "'120.47"
$ws.Range("E30").Value = "'  -2.82%  "
$ws.Range("E31").Value = "'  -4.94%  "
$ws.Range("D32").Value = "'0.09183"
$ws.Range("E32").Value = "'  +0.06%  "
$ws.Range("D33").Value = "'3.603"
$ws.Range("D34").Value = "'5.304"
$ws.Range("E34").Value = "'  -5.38%  "
$ws.Range("D35").Value = "'1.471"
$ws.Range("E35").Value = "'  +6.18%  "
$ws.Range("D36").Value = "'0.02170"
$ws.Range("E36").Value = "'  -4.83%  "
$ws.Range("D37").Value = "'0.05825"
$ws.Range("E37").Value = "'  -4.66%  "
$ws.Range("B38").Value = "'Algorand"
$ws.Range("C38").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.1993"
$ws.Range("E38").Value = "'  -4.66%  "
$ws.Range("B39").Value = "'Aptos"
$ws.Range("C39").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'10.90"
$ws.Range("E39").Value = "'  -7.78%  "
$ws.Range("E40").Value = "'  +0.07%  "
$ws.Range("D41").Value = "'4.699"
$ws.Range("E41").Value = "'  -4.88%  "
$ws.Range("D42").Value = "'0.5954"
$ws.Range("E42").Value = "'  -4.67%  "
$ws.Range("D43").Value = "'1.083"
$ws.Range("E43").Value = "'  -8.05%  "
$ws.Range("D44").Value = "'7.474"
$ws.Range("E44").Value = "'  -3.94%  "
$ws.Range("D45").Value = "'12.78"
$ws.Range("E45").Value = "'  -2.90%  "
$ws.Range("D46").Value = "'3.579"
$ws.Range("E46").Value = "'  -4.22%  "
$ws.Range("D47").Value = "'0.5561"
$ws.Range("E47").Value = "'  -4.88%  "
$ws.Range("D48").Value = "'118.78"
$ws.Range("E48").Value = "'  -3.05%  "
$ws.Range("D50").Value = "'1.105"
$ws.Range("E50").Value = "'  -2.16%  "
$ws.Range("D51").Value = "'0.06631"
$ws.Range("E51").Value = "'  -2.89%  "
